$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 17 (ALC)
$ws.Range("H17").Value = 427.64865
$ws.Range("J17").Value = 427.64865
$ws.Range("L17").Value = 1282.94595
$ws.Range("N17").Value = -1618.94595

# row 41 (ALC)
$ws.Range("H41").Value = 570.1818
$ws.Range("I41").Value = 288.57144
$ws.Range("J41").Value = 1063
$ws.Range("K41").Value = 288.57144
$ws.Range("L41").Value = 1063
$ws.Range("M41").Value = 151.42856
$ws.Range("N41").Value = -1943

# row 86 (ALC)
$ws.Range("H86").Value = 2425.111
$ws.Range("J86").Value = 2921.5
$ws.Range("L86").Value = 2921.5
$ws.Range("N86").Value = -5167.5

# row 88 (ALC)
$ws.Range("H88").Value = 1253399.6
$ws.Range("J88").Value = 1669944
$ws.Range("L88").Value = 1669944
$ws.Range("N88").Value = -1670756

# row 89 (ALC)
$ws.Range("H89").Value = 2425.111
$ws.Range("J89").Value = 2921.5
$ws.Range("L89").Value = 14607.5
$ws.Range("N89").Value = -25839.5

# row 91 (ALC)
$ws.Range("H91").Value = 1253399.6
$ws.Range("J91").Value = 1669944
$ws.Range("L91").Value = 1669944
$ws.Range("N91").Value = -1672752

# row 100 (ALC)
$ws.Range("H100").Value = 3202.0417
$ws.Range("I100").Value = 1857.875
$ws.Range("J100").Value = 5890.375
$ws.Range("K100").Value = 1857.875
$ws.Range("L100").Value = 5890.375
$ws.Range("M100").Value = -1316.875
$ws.Range("N100").Value = -6972.375

# row 137 (ALC)
$ws.Range("H137").Value = 2135.65
$ws.Range("I137").Value = 2491.3845
$ws.Range("J137").Value = 1475
$ws.Range("K137").Value = 7474.1535
$ws.Range("L137").Value = 4425
$ws.Range("M137").Value = -4924.1535
$ws.Range("N137").Value = -9525

# row 138 (ALC)
$ws.Range("H138").Value = 2233.3076
$ws.Range("J138").Value = 2403.8708
$ws.Range("L138").Value = 7211.6124
$ws.Range("N138").Value = -17491.6124

$ws = $wb.Worksheets.Item("ARM")
# row 2 (ARM)
$ws.Range("H2").Value = 3042.9443
$ws.Range("J2").Value = 5470.143
$ws.Range("L2").Value = 5470.143
$ws.Range("N2").Value = -5696.143

# row 32 (ARM)
$ws.Range("H32").Value = 4091.745
$ws.Range("I32").Value = 3558.739
$ws.Range("K32").Value = 3558.739
$ws.Range("M32").Value = -3271.739

# row 61 (ARM)
$ws.Range("H61").Value = 7729.8
$ws.Range("I61").Value = 7500
$ws.Range("J61").Value = 7883
$ws.Range("K61").Value = 7500
$ws.Range("L61").Value = 7883
$ws.Range("M61").Value = -7288
$ws.Range("N61").Value = -8307

# row 116 (ARM)
$ws.Range("H116").Value = 3042.9443
$ws.Range("J116").Value = 5470.143
$ws.Range("L116").Value = 5470.143
$ws.Range("N116").Value = -10058.143

# row 132 (ARM)
$ws.Range("H132").Value = 5395.6313
$ws.Range("I132").Value = 6226.3687
$ws.Range("J132").Value = 3734.158
$ws.Range("K132").Value = 18679.1061
$ws.Range("L132").Value = 11202.474
$ws.Range("M132").Value = -16149.1061
$ws.Range("N132").Value = -16262.474

# row 136 (ARM)
$ws.Range("H136").Value = 7729.8
$ws.Range("I136").Value = 7500
$ws.Range("J136").Value = 7883
$ws.Range("K136").Value = 22500
$ws.Range("L136").Value = 23649
$ws.Range("M136").Value = -19950
$ws.Range("N136").Value = -28749

$ws = $wb.Worksheets.Item("BSM")
# row 3 (BSM)
$ws.Range("H3").Value = 3042.9443
$ws.Range("J3").Value = 5470.143
$ws.Range("L3").Value = 5470.143
$ws.Range("N3").Value = -5698.143

# row 105 (BSM)
$ws.Range("H105").Value = 3101.125
$ws.Range("I105").Value = 2468.3333
$ws.Range("K105").Value = 2468.3333
$ws.Range("M105").Value = -721.3332999999998

$ws = $wb.Worksheets.Item("CRP")
# row 55 (CRP)
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").ClearContents()
$ws.Range("N55").Value = 0

# row 58 (CRP)
$ws.Range("H58").Value = 3013
$ws.Range("I58").Value = 2848.5
$ws.Range("K58").Value = 2848.5
$ws.Range("M58").Value = -2645.5

# row 106 (CRP)
$ws.Range("H106").Value = 100000
$ws.Range("J106").Value = 100000
$ws.Range("L106").Value = 100000
$ws.Range("M106").Value = -102524

# row 118 (CRP)
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").ClearContents()
$ws.Range("N118").Value = 0

# row 132 (CRP)
$ws.Range("H132").Value = 2015.3572
$ws.Range("I132").Value = 1446.7142
$ws.Range("J132").Value = 2584
$ws.Range("K132").Value = 4340.142599999999
$ws.Range("L132").Value = 7752
$ws.Range("M132").Value = -1810.142599999999
$ws.Range("N132").Value = -12812

# row 134 (CRP)
$ws.Range("H134").Value = 2948.0667
$ws.Range("I134").Value = 2833.6667
$ws.Range("J134").Value = 3119.6667
$ws.Range("K134").Value = 8501.000100000001
$ws.Range("L134").Value = 9359.000100000001
$ws.Range("M134").Value = -5966.000100000001
$ws.Range("N134").Value = -14429.0001

# row 136 (CRP)
$ws.Range("H136").Value = 3013
$ws.Range("I136").Value = 2848.5
$ws.Range("K136").Value = 8545.5
$ws.Range("M136").Value = -5995.5

# row 137 (CRP)
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").ClearContents()
$ws.Range("N137").Value = 0

$ws = $wb.Worksheets.Item("CUL")
# row 139 (CUL)
$ws.Range("H139").Value = 3952.3462
$ws.Range("I139").Value = 1058.0714
$ws.Range("J139").Value = 7329
$ws.Range("K139").Value = 3174.2142
$ws.Range("L139").Value = 21987
$ws.Range("M139").Value = 1965.7858
$ws.Range("N139").Value = -32267

# row 140 (CUL)
$ws.Range("H140").Value = 657.4167
$ws.Range("I140").Value = 657.4167
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 1972.2501
$ws.Range("L140").Value = 0
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = 3207.7499

$ws = $wb.Worksheets.Item("GSM")
# row 97 (GSM)
$ws.Range("H97").Value = 745.7
$ws.Range("I97").Value = 626.5
$ws.Range("K97").Value = 626.5
$ws.Range("M97").Value = -130.5

$ws = $wb.Worksheets.Item("LTW")
# row 7 (LTW)
$ws.Range("H7").Value = 31255050
$ws.Range("I7").Value = 5058
$ws.Range("K7").Value = 5058
$ws.Range("M7").Value = -4946

# row 16 (LTW)
$ws.Range("H16").Value = 1646.2084
$ws.Range("I16").Value = 1265.7222
$ws.Range("K16").Value = 1265.7222
$ws.Range("M16").Value = -1095.7222

# row 82 (LTW)
$ws.Range("H82").Value = 3838.05
$ws.Range("J82").Value = 6183.1816
$ws.Range("L82").Value = 6183.1816
$ws.Range("N82").Value = -6905.1816

# row 85 (LTW)
$ws.Range("H85").Value = 3838.05
$ws.Range("J85").Value = 6183.1816
$ws.Range("L85").Value = 6183.1816
$ws.Range("N85").Value = -8679.1816

# row 100 (LTW)
$ws.Range("H100").Value = 3366
$ws.Range("I100").Value = 3007.077
$ws.Range("J100").Value = 4299.2
$ws.Range("K100").Value = 3007.077
$ws.Range("L100").Value = 4299.2
$ws.Range("M100").Value = -2466.077
$ws.Range("N100").Value = -5381.2

# row 122 (LTW)
$ws.Range("H122").Value = 3327.8462
$ws.Range("I122").Value = 3210.9546
$ws.Range("J122").Value = 3479.1177
$ws.Range("K122").Value = 9632.863799999999
$ws.Range("L122").Value = 10437.3531
$ws.Range("M122").Value = -7182.863799999999
$ws.Range("N122").Value = -15337.3531

# row 126 (LTW)
$ws.Range("H126").Value = 31255050
$ws.Range("I126").Value = 5058
$ws.Range("K126").Value = 15174
$ws.Range("M126").Value = -12704

# row 136 (LTW)
$ws.Range("H136").Value = 2852.389
$ws.Range("I136").Value = 1826.6666
$ws.Range("J136").Value = 4288.4
$ws.Range("K136").Value = 5479.9998
$ws.Range("L136").Value = 12865.2
$ws.Range("M136").Value = -2929.9998
$ws.Range("N136").Value = -17965.2

$ws = $wb.Worksheets.Item("WVR")
# row 132 (WVR)
$ws.Range("H132").Value = 4952.1377
$ws.Range("I132").Value = 4852.5
$ws.Range("J132").Value = 5045.1333
$ws.Range("K132").Value = 14557.5
$ws.Range("L132").Value = 15135.3999
$ws.Range("M132").Value = -12027.5
$ws.Range("N132").Value = -20195.3999

# row 141 (WVR)
$ws.Range("H141").Value = 234425.5
$ws.Range("J141").Value = 234425.5
$ws.Range("L141").Value = 234425.5
$ws.Range("N141").Value = -244785.5
